$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update translated text content for rows 2-15 (cells C2:C15)
$ws.Range("C2").Value = 'どんなの能力でも見つけられたら何々を見つけ 開花で5代で始めます じゃあちょっとせ こいかもしれないですけどどんな能力を身に付けられる能力を身につけます'
$ws.Range("C3").Value = 'ドラゴンボールをいつでも作れる能力とかミッション第1回あの人の意見 聞いて 僕は空を飛べたらええなぁ っていうのがあったらいいと思います 自分の過去に戻れるなら僕は透明になれる能力が欲しいな と思います じゃあ僕は見に行ける能力を見つけられたよ じゃあ僕は過去も未来も 両方'
$ws.Range("C4").Value = '過去も未来も両方いける能力が欲しいですね ためなる瞬間移動かなまずまず 瞬間移動と移動時間が一切入らなくなるといえばバイトに行く時間とか学校に通う時間通勤時間 瞬間移動できれば全部なんかも 数秒で終わるてことはそれだけ人生において有意義に使える時間が増えるって事になるからちょっと まあ 普通に便利 力があると生活にもっともっと余裕があるかな'
$ws.Range("C5").Value = '空を飛ぶ お前らんだ理由とした中で多様な意見になっちゃうんですけども アルク より,少しでも早く移動できる空 飛んで移動した方が少しうまくできるかな って思った私の過去に戻りたい子に戻れる能力を得たい理由としましては 昔 楽しかった出来事とかもう一度体験したり その失敗したことをやり直してその今もうそのその失敗を直して その未来を変えるって言うとはできると思うので'
$ws.Range("C6").Value = '透明になりたいっていうのは あの 単純にあの立ちションしたいときとかに止め剤とかでも透明なったら後は女 女子トイレが混んでる時とかに結構入れる'
$ws.Range("C7").Value = 'それぐらいかなと思ったんですけど 違いが分かる かすれば その未来が見える能力を持っていてまーす まーす'
$ws.Range("C8").Value = 'お前人間じゃない 瞬間移動 以外に 過去と未来に活きる能力にもた時点であげてたんですけど 何で出会ってたかって言うの 今 説明するとまずは 加工は自分の何か やり残したとことかをやりたいってのもあるしあと未来に活きる能力は未来に言ったら まあ大葉とかボートレースとかそういうギャンブルの結果が分かるから先に提出の結果を見て現在に戻ってきてでお金を全部入れたらそのぶん 全部返ってくる その能力があるだけで一生お金に困らないでほしい'
$ws.Range("C9").Value = 'その未来に行けたとして未来の自分を見るわけじゃなくその時に見たくない 自分の理想としない未来だった場合になんか怖いしなんか嫌な気持ちになったりしないのかな と思ってます それに対してそれに対しては例えばその未来が確定する前に見た未来が1ヶ月頃だったら まあ でもそうならないように仕向けると言うか そういう風に今できるんじゃないかと未来の自分は今の自分でも理由なんてそうならない間に 同じ意見ですけど手を加えてはいくらでも修正 軌道修正'
$ws.Range("C10").Value = '柿の種でできるんじゃない子なんて ここに居る人達全員がその男をであるって言うのがマジでまひとつで男であるからにはそのを ポリスの性欲ってものあるわけで透明になるっていうのは誰しもが欲しいの 能力ですね はいはい ひとつだけ 懸念てあげるとしたら自分がおじいちゃんとか年を召した時にもそういったことに興味がなくなったりとかしてしまうと労力が無駄になってしまうのかな 生涯現役 なら問題ないですね 年取ってからのこと考えると いらない能力 かもしれん'
$ws.Range("C11").Value = 'しれないそれは確かに 同意する瞬間に感染すると瞬間移動のその あの利便性とかはあると思うよ 三原理って その前に 言ったってやばい かその高速で移動するんじゃなくて自分が減少させてその場所に持っ日市場を形成するって言うんですか あれってどうするんだったら自分がその中に取り込まれるっていうのは アニメ たら よくあるんですけど そういう危険性は別になく 利便性があくまでも そういった問題がない前提の瞬間移動 なので 紅葉 一切していません'
$ws.Range("C12").Value = '紙に瞬間移動という能力 欲しいのだとしたらその問題は解消されているものとして扱っています 俺はもうあの意見を変えます店間移動が欲しいです あの瞬間移動も未来とか過去もいいけど 僕は今を生きるんで深海魚が欲しいな と思います 前回のにそれは決してやましいことをしてすぐに切れるから欲しいとかそういうわけではないです 移動とかあの戦いとかそういう場面でそういう問題から大丈夫'
$ws.Range("C13").Value = '空を飛ぶのも時間かかる ししんかいのノータイムで行けるかな もう一択 北透明になる 結構良かったと透明になる おじさんになってもいるのか っていうのはどういう理由でガッツリ 柄ですよ 僕は'
$ws.Range("C14").Value = '残り3 残り3分 なんで行けると思う'
$ws.Range("C15").Value = '男5人なんて 透明人間 瞬間移動 どちらが高いんですけど 僅差で瞬間移動の勝利ということで ギリギリの僅差で唯一 空を飛ぶ お願いします どちらにつくのか着きます'

# Clear the text for row 16 (C16) while keeping the cell as an empty string cell
$ws.Range("C16").Value = ""

# Remove rows 17-20 (15.wav - 18.wav) entirely
$ws.Rows("17:20").Delete()

Write-Output "done"
